# Updates the cryptocurrency price/volume table on Sheet1 with freshly
# scraped values (GitHub Actions automated refresh).
#
# Columns: A=rank, B=Coin, C=Link, D=Price, E=Volume(1h)
# Price values are kept as plain text (NumberFormat "@") so that values
# such as "29.098.69", "1.000" or "0.00000000124" are preserved exactly
# as scraped instead of being reinterpreted as numbers/dates by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.098.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.835.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6194"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07456"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2929"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07676"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.842.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.014"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6745"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009145"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.908"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.097.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.083.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.216"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1418"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.517"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.505"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05604"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.123"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.136"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.211"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.843"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7421"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.143"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.660"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.775"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01787"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.212.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.408"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9007"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.982.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5096"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4067"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.171"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05830"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.55%  "
